$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.871.84'
$ws.Range("D2").NumberFormat = "General"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.081.29'
$ws.Range("D3").NumberFormat = "General"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.13%  '
$ws.Range("E3").NumberFormat = "General"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.90'
$ws.Range("D5").NumberFormat = "General"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.47%  '
$ws.Range("E5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.56'
$ws.Range("D6").NumberFormat = "General"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("E6").NumberFormat = "General"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -10.00%  '
$ws.Range("E7").NumberFormat = "General"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.359'
$ws.Range("D8").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.68%  '
$ws.Range("E8").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.077.91'
$ws.Range("D10").NumberFormat = "General"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.724'
$ws.Range("D11").NumberFormat = "General"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.07%  '
$ws.Range("E11").NumberFormat = "General"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("E12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("D13").NumberFormat = "General"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.00'
$ws.Range("D14").NumberFormat = "General"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.885.78'
$ws.Range("D15").NumberFormat = "General"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.38'
$ws.Range("D16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.653.74'
$ws.Range("D17").NumberFormat = "General"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.084.98'
$ws.Range("D18").NumberFormat = "General"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.28%  '
$ws.Range("E18").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.81'
$ws.Range("D19").NumberFormat = "General"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("E19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000210'
$ws.Range("D20").NumberFormat = "General"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.81'
$ws.Range("D21").NumberFormat = "General"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.34%  '
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.40'
$ws.Range("D22").NumberFormat = "General"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -9.53%  '
$ws.Range("E22").NumberFormat = "General"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.44'
$ws.Range("D23").NumberFormat = "General"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("E23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.75'
$ws.Range("D24").NumberFormat = "General"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.20%  '
$ws.Range("E24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.58'
$ws.Range("D25").NumberFormat = "General"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.24%  '
$ws.Range("E25").NumberFormat = "General"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Aptos'
$ws.Range("B26").NumberFormat = "General"

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C26").NumberFormat = "General"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.79'
$ws.Range("D26").NumberFormat = "General"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.23%  '
$ws.Range("E26").NumberFormat = "General"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("B27").NumberFormat = "General"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C27").NumberFormat = "General"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '81.51'
$ws.Range("D27").NumberFormat = "General"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -15.81%  '
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.236.43'
$ws.Range("D28").NumberFormat = "General"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("E28").NumberFormat = "General"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E29").NumberFormat = "General"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.13'
$ws.Range("D30").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("E30").NumberFormat = "General"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").NumberFormat = "General"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E31").NumberFormat = "General"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.29%  '
$ws.Range("E32").NumberFormat = "General"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("E33").NumberFormat = "General"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.57'
$ws.Range("D34").NumberFormat = "General"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -10.29%  '
$ws.Range("E34").NumberFormat = "General"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.150'
$ws.Range("D35").NumberFormat = "General"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.65%  '
$ws.Range("E35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.70'
$ws.Range("D36").NumberFormat = "General"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("E36").NumberFormat = "General"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '499.33'
$ws.Range("D37").NumberFormat = "General"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.64%  '
$ws.Range("E37").NumberFormat = "General"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.04'
$ws.Range("D38").NumberFormat = "General"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E38").NumberFormat = "General"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("E39").NumberFormat = "General"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.50%  '
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.60'
$ws.Range("D41").NumberFormat = "General"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +54.95%  '
$ws.Range("E41").NumberFormat = "General"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.34%  '
$ws.Range("E42").NumberFormat = "General"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.10'
$ws.Range("D43").NumberFormat = "General"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("E43").NumberFormat = "General"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.400'
$ws.Range("D45").NumberFormat = "General"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.41%  '
$ws.Range("E45").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.86'
$ws.Range("D46").NumberFormat = "General"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.75%  '
$ws.Range("E46").NumberFormat = "General"

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("B47").NumberFormat = "General"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C47").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.675'
$ws.Range("D47").NumberFormat = "General"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.37%  '
$ws.Range("E47").NumberFormat = "General"

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Monero'
$ws.Range("B48").NumberFormat = "General"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C48").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '149.13'
$ws.Range("D48").NumberFormat = "General"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("E48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.40'
$ws.Range("D49").NumberFormat = "General"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.57%  '
$ws.Range("E49").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.98%  '
$ws.Range("E50").NumberFormat = "General"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").NumberFormat = "General"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.29%  '
$ws.Range("E51").NumberFormat = "General"
